# №10265 от 25.03.2024 https://2eurostore.ru/
# Mark several Portugal 2€ commemorative varieties as owned (0 -> 1 in column F)
# and leave the selection on the last-touched cell (F37).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2€")

$rows = @(20, 22, 23, 27, 28, 29, 30, 31, 32, 33, 34)
foreach ($r in $rows) {
    $ws.Range("F$r").Value = 1
}

$ws.Activate()
$ws.Range("F37").Select()
